$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, count(B), image(C), word(D), category(E)
$rows = @(
    @(2, 34, "house/house011.jpg", "schicken", "house"),
    @(3, 74, "face/face004.jpg", "enden", "face"),
    @(4, 44, "face/face027.jpg", "sieben", "face"),
    @(5, 91, "face/face015.jpg", "nehmen", "face"),
    @(6, 17, "house/house010.jpg", "opfern", "house"),
    @(7, 123, "house/house029.jpg", "hupen", "house"),
    @(8, 92, "face/face013.jpg", "scheitern", "face"),
    @(9, 32, "face/face009.jpg", "füttern", "face"),
    @(10, 29, "house/house019.jpg", "fliegen", "house"),
    @(11, 11, "house/house020.jpg", "schätzen", "house"),
    @(12, 88, "face/face018.jpg", "stechen", "face"),
    @(13, 70, "face/face011.jpg", "töten", "face"),
    @(14, 56, "house/house015.jpg", "bitten", "house"),
    @(15, 38, "face/face022.jpg", "schmecken", "face"),
    @(16, 105, "face/face025.jpg", "liefern", "face"),
    @(17, 35, "house/house026.jpg", "krachen", "house"),
    @(18, 113, "house/house028.jpg", "biegen", "house"),
    @(19, 93, "house/house003.jpg", "ehren", "house"),
    @(20, 100, "face/face030.jpg", "kaufen", "face"),
    @(21, 94, "face/face029.jpg", "drehen", "face"),
    @(22, 49, "house/house022.jpg", "segeln", "house"),
    @(23, 14, "face/face026.jpg", "füllen", "face"),
    @(24, 15, "house/house031.jpg", "tagen", "house"),
    @(25, 22, "house/house016.jpg", "raten", "house"),
    @(26, 13, "face/face020.jpg", "spielen", "face"),
    @(27, 2, "house/house009.jpg", "formen", "house"),
    @(28, 90, "house/house018.jpg", "starten", "house"),
    @(29, 30, "face/face031.jpg", "regnen", "face"),
    @(30, 18, "face/face010.jpg", "wenden", "face"),
    @(31, 43, "house/house021.jpg", "posten", "house"),
    @(32, 54, "house/house000.jpg", "loben", "house"),
    @(33, 106, "face/face023.jpg", "rücken", "face")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}
